# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" column (E) on Hoja1 lists 22 billing periods for rows
# 16-37, originally ascending from 1806 (Jun-2018) through 2003 (Mar-2020).
# This edit removes that previous list of periods and replaces it with the
# same 22 periods in reverse (descending) order, 2003 down to 1806 -- i.e.
# the old "Estados de Cuenta" periods are dropped and re-added newest-first.
# All other data in each row (Tipo Doc, N Doc, Nombre, Valor Mora, Salario
# Basico) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New (descending) period order replacing the old ascending one.
$periodos = @(
    "2003", "2002", "2001",
    "1912", "1911", "1910", "1909", "1908", "1907", "1906", "1905", "1904", "1903", "1902", "1901",
    "1812", "1811", "1810", "1809", "1808", "1807", "1806"
)

$firstRow = 16
$lastRow = 37

# Drop the previous periods first ...
$ws.Range("E" + $firstRow + ":E" + $lastRow).ClearContents()

# ... then add the new ones in the desired (reversed) order.
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("E" + $row).Value = $periodos[$i]
}
